# T1063_Contacts_AddMultipleContactsToAnExistingCompanyAndValidateSaveAndNewButton.xlsx
# "Contacts Final - 14 Oct 2024"
#
# The "PhysicalOffice" / "DC" column (column O) on the Contact sheet is no
# longer needed, so it is removed entirely. Removing the whole column shifts
# every column to its right (Title, Department, LineOfBusiness, PinCode) one
# position to the left, and also drops the two now-unused shared strings
# ("PhysicalOffice" and "DC") from the workbook once Excel re-saves it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Delete the whole "PhysicalOffice"/"DC" column (column O), shifting the
# remaining columns (Title, Department, LineOfBusiness, PinCode) left.
$ws.Columns("O:O").Delete()

# Leave the sheet scrolled/selected the way the author left it when saving.
$ws.Activate()
$ws.Range("H12").Select()
